# Gantt chart update: append two more data rows (23-Jan / 24-Jan / 25-Jan cycle
# -> actually two new rows for 44950 and 44951) to "Main Board", push the
# trailing thick-border spacer rows down, extend the chart series ranges to
# match, and shift the chart's anchor down by the same two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Board")

# --- 1. Make room: insert two new rows above the current spacer rows (7 & 8)
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# --- 2. Populate the two new data rows
$ws.Range("A7").Value = 44950
$ws.Range("B7").Value = 37
$ws.Range("C7").Value = 12

$ws.Range("A8").Value = 44951
$ws.Range("B8").Value = 39
$ws.Range("C8").Value = 12

# --- 3. Match the date formatting/style used by the rest of column A
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- 4. Move the selection the way the saved file records it
$ws.Range("D8").Select()

# --- 5. Extend the two chart series so they cover the new rows too
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$sClosed = $chart.SeriesCollection().Item(1)
$sOpen = $chart.SeriesCollection().Item(2)

$sClosed.Formula = "=SERIES('Main Board'!`$C`$1,'Main Board'!`$A`$2:`$A`$9,'Main Board'!`$C`$2:`$C`$9,1)"
$sOpen.Formula = "=SERIES('Main Board'!`$B`$1,'Main Board'!`$A`$2:`$A`$9,'Main Board'!`$B`$2:`$B`$9,2)"

# --- 6. The chart is cell-anchored; nudge it down by the height of the two
#        rows we just inserted so its anchor tracks the data the way it did
#        before (twoCellAnchor "move and size with cells").
$co.Top = $co.Top + 30.0
